$d = $word.ActiveDocument

# 1. Insert the new "normalized to baseline (A)" sentence in place of the old
#    "A) and baseline blood glucose (B) were measured in mice" lead-in, and
#    splice in the new insulin-injection / cohort sentence before the old
#    "and baseline blood glucose (" text.
$d.Content.Find.Execute(
    "A) and baseline blood glucose (B) were measured in mice",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "normalized to baseline (A). Insulin was given via i.p. injection at a concentration of 2.5 U/kg following five weeks of dexamethasone (NCD n=12; HFD n=12) or vehicle (NCD n=12; HFD n=12) treatment and 17 weeks of diet. ITT (B) and baseline blood glucose (C) were measured in another cohort of mice",
    2) | Out-Null

# 2. Relabel the fat/lean mass panel letters C -> D and D -> E.
$d.Content.Find.Execute(
    "Fat (C) and lean mass (D) were measured weekly",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Fat (D) and lean mass (E) were measured weekly",
    2) | Out-Null

# 3. Relabel the blood glucose / insulin clearance panel letters E -> F, F -> G.
$d.Content.Find.Execute(
    "blood glucose levels (E), insulin clearance rates (F),",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "blood glucose levels (F), insulin clearance rates (G),",
    2) | Out-Null

# 4. Relabel the tissue-uptake panel letters G -> H, H -> I. (Anchor begins
#    right after the "iWAT" run/proofErr-end boundary so the spell-check
#    proofErr start/end pair around "iWAT" stays intact.)
$d.Content.Find.Execute(
    " (G) and heart and brown adipose tissue (H) were measured",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " (H) and heart and brown adipose tissue (I) were measured",
    2) | Out-Null
